# Update "想去人数" (interested-people count) figures in column F across
# the 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) sheets to
# reflect the freshly scraped gh-pages output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibit.Range("F2").Value  = 52
$wsExhibit.Range("F5").Value  = 1073
$wsExhibit.Range("F8").Value  = 582
$wsExhibit.Range("F11").Value = 1404
$wsExhibit.Range("F12").Value = 3044
$wsExhibit.Range("F13").Value = 548
$wsExhibit.Range("F14").Value = 1706
$wsExhibit.Range("F16").Value = 828
$wsExhibit.Range("F18").Value = 1435
$wsExhibit.Range("F23").Value = 425
$wsExhibit.Range("F24").Value = 50
$wsExhibit.Range("F25").Value = 4110
$wsExhibit.Range("F26").Value = 721
$wsExhibit.Range("F28").Value = 1602
$wsExhibit.Range("F29").Value = 7
$wsExhibit.Range("F30").Value = 71

# 演出 (sheet2)
$wsShow.Range("F3").Value = 175
$wsShow.Range("F6").Value = 58
$wsShow.Range("F9").Value = 42

# 全部类型 (sheet4)
$wsAll.Range("F2").Value  = 52
$wsAll.Range("F7").Value  = 175
$wsAll.Range("F10").Value = 58
$wsAll.Range("F14").Value = 42
$wsAll.Range("F16").Value = 1073
$wsAll.Range("F19").Value = 582
$wsAll.Range("F22").Value = 1404
$wsAll.Range("F23").Value = 3044
$wsAll.Range("F24").Value = 548
$wsAll.Range("F25").Value = 1706
$wsAll.Range("F27").Value = 828
$wsAll.Range("F29").Value = 1435
$wsAll.Range("F36").Value = 425
$wsAll.Range("F37").Value = 50
$wsAll.Range("F38").Value = 4110
$wsAll.Range("F39").Value = 721
$wsAll.Range("F41").Value = 1602
$wsAll.Range("F44").Value = 7
$wsAll.Range("F45").Value = 71
